$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.000.90"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.572.48"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.04"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.03"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "2.572.82"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.159"
$ws.Range("E12").Value = "  +11.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "3.027.47"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "59.043.74"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.37"
$ws.Range("E16").Value = "  +7.14%  "
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("D18").Value = "2.579.42"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.04"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.21"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.24"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.40"
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.457"
$ws.Range("E25").Value = "  +6.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.24"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  +2.62%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.44"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.97"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.870"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.16"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.50"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.67"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "292.53"
$ws.Range("E42").Value = "  +3.71%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0977"
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.591"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.18"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.27"
$ws.Range("E49").Value = "  +5.52%  "
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("D51").Value = "1.945.55"
